$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 28-31: mark "Hoan thanh" (E) as 100% and set the assignee (F) to the
# new combined name "Bùi, Kiều". Row 28 previously had this note in column G
# ("Bùi") - remove it since the note moves to column F with the new text.
$ws.Range("E28").Value = 1
$ws.Range("G28").Clear()
$ws.Range("F28").Value = "Bùi, Kiều"

$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Bùi, Kiều"

$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Bùi, Kiều"

$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "Bùi, Kiều"

# Rows 44-45: mark "Hoan thanh" (E) as 100%.
$ws.Range("E44").Value = 1
$ws.Range("E45").Value = 1

# Update the view: scroll to A13 and select G30.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G30").Select()
